$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on the existing two sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row styling to match the other sheets (bold, centered, thin border)
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Date-format the "ds" column like the other sheets
$wsForecast.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$rows = @(
    @(44983.99999999999, 51, 11.03207632652852, 94.1719491726483),
    @(44990.99999999999, 49, 8.199735782023552, 95.28206728183594),
    @(45067.99999999999, 19, -24.81250189460252, 64.48516742824351),
    @(45074.99999999999, 16, -23.22435315505464, 59.45855693974851),
    @(45081.99999999999, 13, -27.57830308999292, 53.62647216231088),
    @(45088.99999999999, 11, -33.57444207668791, 52.47818999558696),
    @(45095.99999999999, 8, -36.57250790812871, 51.77484288567084),
    @(45102.99999999999, 5, -37.10483706344058, 48.55072764562655),
    @(45109.99999999999, 3, -41.49915585866663, 45.71657761770725),
    @(45116.99999999999, 0, -44.62799247834445, 40.52451740663066),
    @(45123.99999999999, 0, -48.52461650437161, 39.51449528337304),
    @(45130.99999999999, 0, -47.02084959027239, 34.6664422549397)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wsForecast.Range("A1").Select()
